$wb = $excel.ActiveWorkbook

# --- 1. Rename the first sheet: "default" -> "nauwkeurigheid" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "nauwkeurigheid"

$ws2 = $wb.Worksheets.Item(2)

# --- 2. Add a new worksheet "snelheid" after the existing sheets ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$newSheet.Name = "snelheid"

# --- 3. Populate the data for the 2nd test report ---
$newSheet.Range("A1").Value = "Test-ID"
$newSheet.Range("B1").Value = "Computer1:"
$newSheet.Range("C1").Value = "Computer2:"
$newSheet.Range("D1").Value = "Computer3:"

$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 2057
$newSheet.Range("C2").Value = 5694

$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = 2538
$newSheet.Range("C3").Value = 13090

$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = 10911
$newSheet.Range("C4").Value = 29199

# --- 4. Turn the range into a table ("Tabel12") ---
$listObj = $newSheet.ListObjects.Add(1, $newSheet.Range("A1:D4"), [System.Reflection.Missing]::Value, 1)
$listObj.Name = "Tabel12"
$listObj.TableStyle = "TableStyleMedium8"

# --- 5. Show the totals row with the custom ratio formulas ---
$listObj.ShowTotals = $true

$newSheet.Range("A5").Value = "test3/test1"
$newSheet.Range("B5").Formula = "=B4/B2"
$newSheet.Range("C5").Formula = "=C4/C2"

$col1 = $listObj.ListColumns.Item(1)
$col1.TotalsRowLabel = "test3/test1"

$col2 = $listObj.ListColumns.Item(2)
$col2.TotalsRowFormula = "=B4/B2"

$col3 = $listObj.ListColumns.Item(3)
$col3.TotalsRowFormula = "=C4/C2"

# --- 6. Page setup for the new sheet ---
$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# --- 7. Restore/record selections on each sheet ---
$ws1.Activate() | Out-Null
$ws1.Range("C31").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B5").Select() | Out-Null

# --- 8. Make "snelheid" the active/visible tab, selecting D5 ---
$newSheet.Activate() | Out-Null
$newSheet.Range("D5").Select() | Out-Null
